$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.24'
$ws.Range("G2").Value = '''17'
$ws.Range("G3").Value = '''17'
$ws.Range("D4").Value = '''5.450'
$ws.Range("G4").Value = '''17'
$ws.Range("D5").Value = '''0.05615'
$ws.Range("G5").Value = '''17'
$ws.Range("D6").Value = '''6.458'
$ws.Range("G6").Value = '''17'
$ws.Range("D7").Value = '''0.8044'
$ws.Range("G7").Value = '''17'
$ws.Range("G8").Value = '''17'
$ws.Range("D9").Value = '''0.1424'
$ws.Range("G9").Value = '''17'
$ws.Range("D10").Value = '''0.07323'
$ws.Range("G10").Value = '''17'
$ws.Range("G11").Value = '''17'
$ws.Range("D12").Value = '''0.02932'
$ws.Range("G12").Value = '''17'
$ws.Range("D13").Value = '''0.09260'
$ws.Range("G13").Value = '''17'
$ws.Range("D14").Value = '''0.001666'
$ws.Range("G14").Value = '''17'
$ws.Range("D15").Value = '''3.214'
$ws.Range("G15").Value = '''17'
$ws.Range("D16").Value = '''0.04741'
$ws.Range("G16").Value = '''17'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.0005822'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("G17").Value = '''17'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '''0.006372'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("G18").Value = '''17'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '''0.005070'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("G19").Value = '''17'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = '''0.001054'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("G20").Value = '''17'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = '''0.0001503'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("G21").Value = '''17'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = '''3.989'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("G22").Value = '''17'
$ws.Range("B23").Value = 'GateToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D23").Value = '''3.382'
$ws.Range("E23").Value = '22GateTokenGT'
$ws.Range("G23").Value = '''17'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.122'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("G24").Value = '''17'
$ws.Range("G25").Value = '''17'
$ws.Range("E26").Value = '25ProBitTokenPROBBestin24h'
$ws.Range("G26").Value = '''17'
$ws.Range("G27").Value = '''17'
$ws.Range("G28").Value = '''17'
$ws.Range("G29").Value = '''17'
$ws.Range("G30").Value = '''17'
$ws.Range("G31").Value = '''17'
$ws.Range("G32").Value = '''17'
$ws.Range("G33").Value = '''17'
$ws.Range("G34").Value = '''17'
$ws.Range("G35").Value = '''17'
$ws.Range("G36").Value = '''17'
$ws.Range("G37").Value = '''17'
$ws.Range("G38").Value = '''17'
$ws.Range("G39").Value = '''17'
$ws.Range("D40").Value = '''0.04166'
$ws.Range("G40").Value = '''17'
$ws.Range("D41").Value = '''0.006906'
$ws.Range("G41").Value = '''17'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1040'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '''17'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002977'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").Value = '''17'
$ws.Range("D44").Value = '''0.008540'
$ws.Range("G44").Value = '''17'
$ws.Range("D45").Value = '''0.00005659'
$ws.Range("G45").Value = '''17'
$ws.Range("D46").Value = '''0.00000000751'
$ws.Range("G46").Value = '''17'
$ws.Range("D47").Value = '''0.6813'
$ws.Range("G47").Value = '''17'
$ws.Range("D48").Value = '''0.01518'
$ws.Range("G48").Value = '''17'
$ws.Range("D49").Value = '''0.00002104'
$ws.Range("G49").Value = '''17'
$ws.Range("G50").Value = '''17'
$ws.Range("G51").Value = '''17'
